$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-5, 12, 14: only the "Förändrad" (C) column advances to the new date serial
$ws.Range("C2").Value = 46074
$ws.Range("C3").Value = 46074
$ws.Range("C4").Value = 46074
$ws.Range("C5").Value = 46074

# Row 6
$ws.Range("A6").Value = "A 2593-2024"
$ws.Range("B6").Value = 45313.69204861111
$ws.Range("C6").Value = 46074
$ws.Range("G6").Value = 2.3

# Row 7
$ws.Range("A7").Value = "A 12651-2022"
$ws.Range("B7").Value = 44641
$ws.Range("C7").Value = 46074
$ws.Range("G7").Value = 3.2

# Row 8
$ws.Range("A8").Value = "A 5792-2024"
$ws.Range("B8").Value = 45335
$ws.Range("C8").Value = 46074
$ws.Range("G8").Value = 5.6

# Row 9
$ws.Range("A9").Value = "A 13651-2023"
$ws.Range("B9").Value = 45006
$ws.Range("C9").Value = 46074
$ws.Range("G9").Value = 2.2

# Row 10
$ws.Range("A10").Value = "A 8194-2025"
$ws.Range("B10").Value = 45708
$ws.Range("C10").Value = 46074
$ws.Range("G10").Value = 1.9

# Row 11
$ws.Range("A11").Value = "A 50997-2025"
$ws.Range("B11").Value = 45946
$ws.Range("C11").Value = 46074
$ws.Range("G11").Value = 1.5

# Row 12
$ws.Range("C12").Value = 46074

# Row 13
$ws.Range("A13").Value = "A 7827-2026"
$ws.Range("B13").Value = 46062.63958333333
$ws.Range("C13").Value = 46074
$ws.Range("G13").Value = 2.1

# Row 14
$ws.Range("C14").Value = 46074

# Row 15
$ws.Range("A15").Value = "A 28288-2023"
$ws.Range("B15").Value = 45099.6349537037
$ws.Range("C15").Value = 46074
$ws.Range("G15").Value = 0.5

# Row 16
$ws.Range("A16").Value = "A 7333-2025"
$ws.Range("B16").Value = 45703.35899305555
$ws.Range("C16").Value = 46074
$ws.Range("G16").Value = 0.9
